# Correction in SA algorithm and 746 logs
# Updates the "Fitness" column (C) for run_22 log data to reflect the
# corrected simulated-annealing algorithm output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Generations 0-16  (rows 2-18)  -> Fitness 8199
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 3).Value = 8199
}

# Generations 17-53 (rows 19-55) -> Fitness 8150
for ($r = 19; $r -le 55; $r++) {
    $ws.Cells.Item($r, 3).Value = 8150
}

# Generations 54-78 (rows 56-80) -> Fitness 7811
for ($r = 56; $r -le 80; $r++) {
    $ws.Cells.Item($r, 3).Value = 7811
}

# Generations 79-84 (rows 81-86) -> Fitness 7622
for ($r = 81; $r -le 86; $r++) {
    $ws.Cells.Item($r, 3).Value = 7622
}

# Generations 85-250 (rows 87-252) -> Fitness 7573
for ($r = 87; $r -le 252; $r++) {
    $ws.Cells.Item($r, 3).Value = 7573
}
